# Update Name of Algo
# Apply corrected/updated imputed values produced by the RandomForest
# algorithm run to the result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.1742
$ws.Range("C3").Value = -12.1068
$ws.Range("D8").Value = -9.193599999999998
$ws.Range("D11").Value = -7.645299999999999
$ws.Range("A12").Value = -21.6323
$ws.Range("B14").Value = 5.6001
$ws.Range("D14").Value = -7.5503
$ws.Range("D15").Value = -8.412999999999997
$ws.Range("B26").Value = 3.620900000000006
$ws.Range("C30").Value = -13.44209999999999
$ws.Range("B31").Value = 4.733400000000004
$ws.Range("A32").Value = -21.32720000000001
$ws.Range("B35").Value = 9.234000000000004
$ws.Range("A36").Value = -19.8977
$ws.Range("D36").Value = -7.461500000000004
$ws.Range("B37").Value = 8.726399999999998
$ws.Range("A38").Value = -19.7063
$ws.Range("C44").Value = -14.09569999999999
$ws.Range("B45").Value = 5.063900000000002
$ws.Range("A46").Value = -21.80370000000001
$ws.Range("A54").Value = -21.7984
$ws.Range("A55").Value = -21.95809999999999
$ws.Range("B57").Value = 4.931499999999994
$ws.Range("C58").Value = -12.62599999999999
$ws.Range("D64").Value = -7.515399999999996
$ws.Range("A67").Value = -21.50389999999998
$ws.Range("A69").Value = -21.60499999999997
$ws.Range("A72").Value = -21.8734
$ws.Range("C84").Value = -13.84979999999999
$ws.Range("C89").Value = -10.9049
$ws.Range("D89").Value = -5.9533
$ws.Range("A91").Value = -21.52730000000001
$ws.Range("C91").Value = -10.9768
$ws.Range("C92").Value = -11.4892
$ws.Range("A99").Value = -20.1282
$ws.Range("B100").Value = 5.376899999999996
$ws.Range("B102").Value = 8.194600000000001
$ws.Range("C102").Value = -12.66970000000001
